# Update "想去人数" (column F) counts on all sheets to the refreshed
# snapshot values (gh-pages data regeneration @ 456a3b4).
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local life)
# Sheet 4 = 全部类型 (All types, union of sheets 1-3)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7784
$ws.Range("F3").Value = 104
$ws.Range("F4").Value = 81
$ws.Range("F5").Value = 8871
$ws.Range("F8").Value = 642
$ws.Range("F10").Value = 138
$ws.Range("F11").Value = 439
$ws.Range("F12").Value = 777
$ws.Range("F13").Value = 42
$ws.Range("F14").Value = 76
$ws.Range("F15").Value = 317
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 265
$ws.Range("F19").Value = 396
$ws.Range("F20").Value = 150
$ws.Range("F21").Value = 1087
$ws.Range("F23").Value = 625
$ws.Range("F24").Value = 2215
$ws.Range("F25").Value = 738
$ws.Range("F26").Value = 53
$ws.Range("F29").Value = 615

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 328
$ws.Range("F9").Value = 140

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 457

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 457
$ws.Range("F3").Value = 7785
$ws.Range("F4").Value = 104
$ws.Range("F5").Value = 81
$ws.Range("F7").Value = 8872
$ws.Range("F10").Value = 642
$ws.Range("F13").Value = 138
$ws.Range("F14").Value = 439
$ws.Range("F15").Value = 328
$ws.Range("F18").Value = 777
$ws.Range("F19").Value = 42
$ws.Range("F20").Value = 76
$ws.Range("F21").Value = 317
$ws.Range("F23").Value = 18
$ws.Range("F25").Value = 140
$ws.Range("F27").Value = 265
$ws.Range("F29").Value = 396
$ws.Range("F30").Value = 150
$ws.Range("F31").Value = 1087
$ws.Range("F33").Value = 625
$ws.Range("F34").Value = 2215
$ws.Range("F35").Value = 738
$ws.Range("F36").Value = 53
$ws.Range("F40").Value = 615
